# Update NATMI Fgf1-Cspg4 LR-pair stats with recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.7304773333333333
$ws.Range("H2").Value = 2.191432
$ws.Range("I2").Value = 0.03163269997405359
$ws.Range("J2").Value = 0.03163269997405359
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.366995666666667
$ws.Range("N2").Value = 4.100987
$ws.Range("O2").Value = 0.02653821474268573
$ws.Range("P2").Value = 0.02653821474268573
$ws.Range("Q2").Value = 0.9985593492648889
$ws.Range("R2").Value = 8.987034143383999
$ws.Range("S2").Value = 0.0008394753848023838
$ws.Range("T2").Value = 0.0008394753848023835
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.7304773333333333
$ws.Range("H3").Value = 2.191432
$ws.Range("I3").Value = 0.03163269997405359
$ws.Range("J3").Value = 0.03163269997405359
$ws.Range("O3").Value = 0.2493665720274216
$ws.Range("P3").Value = 0.2493665720274215
$ws.Range("Q3").Value = 9.382971850461331
$ws.Range("R3").Value = 84.44674665415198
$ws.Range("S3").Value = 0.007888137956501652
$ws.Range("T3").Value = 0.007888137956501652
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.7304773333333333
$ws.Range("H4").Value = 2.191432
$ws.Range("I4").Value = 0.03163269997405359
$ws.Range("J4").Value = 0.03163269997405359
$ws.Range("M4").Value = 37.298478
$ws.Range("N4").Value = 111.895434
$ws.Range("O4").Value = 0.7240952132298927
$ws.Range("P4").Value = 0.7240952132298926
$ws.Range("Q4").Value = 27.245692746832
$ws.Range("R4").Value = 245.211234721488
$ws.Range("S4").Value = 0.02290508663274956
$ws.Range("T4").Value = 0.02290508663274956
$ws.Range("I5").Value = 0.4074771110502447
$ws.Range("J5").Value = 0.4074771110502448
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.366995666666667
$ws.Range("N5").Value = 4.100987
$ws.Range("O5").Value = 0.02653821474268573
$ws.Range("P5").Value = 0.02653821474268573
$ws.Range("Q5").Value = 12.86295760982833
$ws.Range("R5").Value = 115.766618488455
$ws.Range("S5").Value = 0.0108137150757806
$ws.Range("T5").Value = 0.0108137150757806
$ws.Range("I6").Value = 0.4074771110502447
$ws.Range("J6").Value = 0.4074771110502448
$ws.Range("O6").Value = 0.2493665720274216
$ws.Range("P6").Value = 0.2493665720274215
$ws.Range("S6").Value = 0.1016111703622365
$ws.Range("T6").Value = 0.1016111703622365
$ws.Range("I7").Value = 0.4074771110502447
$ws.Range("J7").Value = 0.4074771110502448
$ws.Range("M7").Value = 37.298478
$ws.Range("N7").Value = 111.895434
$ws.Range("O7").Value = 0.7240952132298927
$ws.Range("P7").Value = 0.7240952132298926
$ws.Range("Q7").Value = 350.9658100050899
$ws.Range("R7").Value = 3158.69229004581
$ws.Range("S7").Value = 0.2950522256122276
$ws.Range("T7").Value = 0.2950522256122276
$ws.Range("G8").Value = 12.95234266666667
$ws.Range("H8").Value = 38.857028
$ws.Range("I8").Value = 0.5608901889757016
$ws.Range("J8").Value = 0.5608901889757018
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.366995666666667
$ws.Range("N8").Value = 4.100987
$ws.Range("O8").Value = 0.02653821474268573
$ws.Range("P8").Value = 0.02653821474268573
$ws.Range("Q8").Value = 17.70579629851511
$ws.Range("R8").Value = 159.352166686636
$ws.Range("S8").Value = 0.01488502428210275
$ws.Range("T8").Value = 0.01488502428210275
$ws.Range("G9").Value = 12.95234266666667
$ws.Range("H9").Value = 38.857028
$ws.Range("I9").Value = 0.5608901889757016
$ws.Range("J9").Value = 0.5608901889757018
$ws.Range("O9").Value = 0.2493665720274216
$ws.Range("P9").Value = 0.2493665720274215
$ws.Range("Q9").Value = 166.3726731728786
$ws.Range("R9").Value = 1497.354058555908
$ws.Range("S9").Value = 0.1398672637086834
$ws.Range("T9").Value = 0.1398672637086834
$ws.Range("G10").Value = 12.95234266666667
$ws.Range("H10").Value = 38.857028
$ws.Range("I10").Value = 0.5608901889757016
$ws.Range("J10").Value = 0.5608901889757018
$ws.Range("M10").Value = 37.298478
$ws.Range("N10").Value = 111.895434
$ws.Range("O10").Value = 0.7240952132298927
$ws.Range("P10").Value = 0.7240952132298926
$ws.Range("Q10").Value = 483.1026680011279
$ws.Range("R10").Value = 4347.924012010152
$ws.Range("S10").Value = 0.4061379009849155
$ws.Range("T10").Value = 0.4061379009849155
